# =====================================================================
# Edit script: converts the "merged report" layout into a flat table.
#  - rename sheet Sheet1 -> data
#  - set explicit column widths for A:E
#  - unmerge every merged cell in column A, filling the blank cells
#    that the merge used to cover with the same (repeated) value
#  - re-point the "general"/left-aligned column-A style so the filled
#    rows read with plain (unset) horizontal alignment + no bottom rule
#  - normalise row heights for the rows whose wrapped text grew/shrank
#  - move the active selection to G6
# =====================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- constants (mirrors Microsoft.Office.Interop.Excel enums) -------
$xlContinuous      = 1
$xlLineStyleNone    = -4142
$xlDash             = -4115
$xlThin             = 2
$xlHAlignGeneral    = 1
$xlVAlignCenter     = -4108
$xlEdgeLeft         = 7
$xlEdgeTop          = 8
$xlEdgeRight        = 10
$xlEdgeBottom       = 9

# ---- 0. sheet name ----------------------------------------------------
$ws.Name = "data"

# ---- 1. column widths --------------------------------------------------
$ws.Range("A:A").ColumnWidth = 42.2640625
$ws.Range("B:B").ColumnWidth = 5.93203125
$ws.Range("C:C").ColumnWidth = 10.2640625
$ws.Range("D:D").ColumnWidth = 18.93203125
$ws.Range("E:E").ColumnWidth = 26.6

# ---- helper: unmerge a column-A group & repeat its value down the rows
function Unmerge-Fill([string]$rangeAddr) {
    $rng = $ws.Range($rangeAddr)
    $val = $rng.Cells.Item(1, 1).Value2
    $rng.UnMerge()
    $rng.Value2 = $val
}

# ---- 2. unmerge + fill every grouped block in column A ----------------
Unmerge-Fill "A2:A3"
Unmerge-Fill "A5:A6"
Unmerge-Fill "A7:A11"
Unmerge-Fill "A13:A14"
Unmerge-Fill "A15:A18"
Unmerge-Fill "A19:A21"
Unmerge-Fill "A22:A25"

# ---- 3. restyle column A (now a plain repeated-value column) ----------
# every row except the stand-alone "big" rows (4, 12, 26) switches to:
#   border: left/right/top continuous thin, bottom none
#   alignment: horizontal = general (cleared), vertical = center, wrap text
$colARows = @(2,3,5,6,7,8,9,10,11,13,14,15,16,17,18,19,20,21,22,23,24,25)
foreach ($r in $colARows) {
    $c = $ws.Cells.Item($r, 1)
    $c.HorizontalAlignment = $xlHAlignGeneral
    $c.VerticalAlignment   = $xlVAlignCenter
    $c.WrapText            = $true
    $c.Borders.Item($xlEdgeLeft).LineStyle   = $xlContinuous
    $c.Borders.Item($xlEdgeLeft).Weight      = $xlThin
    $c.Borders.Item($xlEdgeRight).LineStyle  = $xlContinuous
    $c.Borders.Item($xlEdgeRight).Weight     = $xlThin
    $c.Borders.Item($xlEdgeTop).LineStyle    = $xlContinuous
    $c.Borders.Item($xlEdgeTop).Weight       = $xlThin
    $c.Borders.Item($xlEdgeBottom).LineStyle = $xlLineStyleNone
}

# ---- 4. row heights -----------------------------------------------------
$ws.Rows.Item(2).RowHeight  = 26
$ws.Rows.Item(4).RowHeight  = 65
$ws.Rows.Item(12).RowHeight = 39
$ws.Rows.Item(16).RowHeight = 26
$ws.Rows.Item(17).RowHeight = 26
$ws.Rows.Item(19).RowHeight = 26
$ws.Rows.Item(20).RowHeight = 39
$ws.Rows.Item(21).RowHeight = 39
$ws.Rows.Item(23).RowHeight = 26
$ws.Rows.Item(24).RowHeight = 26
$ws.Rows.Item(26).RowHeight = 52

# ---- 5. selection --------------------------------------------------------
$ws.Range("G6").Select()
